$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Extend the existing header style (bold + border + centered, style index "1"
#        in the original file) to the new cells that need it, BEFORE we touch A1's
#        own formatting. This keeps the style table free of new/duplicate xf records
#        because we are reusing a style that already exists via copy/paste-format
#        instead of re-building it through property assignment (which would create
#        new xf entries for every incremental change).
$ws.Range("A1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:A8").PasteSpecial(-4122)   # xlPasteFormats

# --- 2) Reset cells that must end up with the *default* (unstyled) look: the old
#        A1 slot (dropped entirely) and the whole inner data block B2:N8. Z100 was
#        never touched, so its format is the workbook's default style (index 0) -
#        copying it over is a clean way to strip formatting without creating new
#        style records.
$ws.Range("Z100").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B2:N8").PasteSpecial(-4122)

# --- 3) Write the new values -------------------------------------------------
# Header row: plain running numbers 0..12 across B1:N1 (no more shared-string
# "Estratégia"/"Related"/... labels; the position index of the dataframe columns
# itself).
$headers = 0..12
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Blank out the old A1 content (header style already reset above).
$ws.Range("A1").ClearContents()

# Row index column (A2:A8) = 0..6
$rowIndex = 0..6
for ($i = 0; $i -lt $rowIndex.Count; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $rowIndex[$i]
}

# Strategy names (B2:B8) - new shared strings
$names = @("Search", "Search + SB", "GS + SB", "Scopus + SB", "Scopus + BS // FS", "Scopus + BS + FS", "Scopus + FS + BS")
for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $names[$i]
}

# Data matrix, columns C..N (3..14), rows 2..8
$data = @(
    @(11, 480, 11, 0, 0, 0, 0.02291666666666667, 0.7857142857142857, 0.04453441295546558, 0.02291666666666667, 0.7857142857142857, 0.04453441295546558),
    @(14, 932, 14, 0, 0, 0, 0.01502145922746781, 1, 0.02959830866807611, 0.01502145922746781, 1, 0.02959830866807611),
    @(11, 442, 11, 0, 0, 0, 0.0248868778280543, 0.7857142857142857, 0.04824561403508772, 0.0248868778280543, 0.7857142857142857, 0.04824561403508772),
    @(11, 502, 11, 0, 0, 0, 0.02191235059760956, 0.7857142857142857, 0.04263565891472869, 0.02191235059760956, 0.7857142857142857, 0.04263565891472869),
    @(9, 242, 9, 0, 0, 0, 0.0371900826446281, 0.6428571428571429, 0.0703125, 0.0371900826446281, 0.6428571428571429, 0.0703125),
    @(11, 424, 11, 0, 0, 0, 0.0259433962264151, 0.7857142857142857, 0.0502283105022831, 0.0259433962264151, 0.7857142857142857, 0.0502283105022831),
    @(9, 275, 9, 0, 0, 0, 0.03272727272727273, 0.6428571428571429, 0.06228373702422146, 0.03272727272727273, 0.6428571428571429, 0.06228373702422146)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item(2 + $r, 3 + $c).Value = $rowVals[$c]
    }
}
